$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"
